$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in remaining values on row 33 (C:F) ---
# These need to be stored as TEXT (not numbers), matching the sheet's
# existing convention for this block of rows. Forcing the NumberFormat to
# "@" (Text) before assigning the value makes the engine store it as a
# string instead of auto-converting it to a number; ClearFormats()
# afterwards removes the now-unneeded explicit formatting so the cell is
# left with the default style, same as its neighbours.
$row33 = @{
    "C33" = "5"
    "D33" = "6"
    "E33" = "4"
    "F33" = "-2"
}
foreach ($addr in $row33.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $row33[$addr]
    $cell.ClearFormats()
}

# --- Add new row 34 ("tela de atuação user") ---
$row34 = @{
    "A34" = "hoje"
    "B34" = "nome"
    "C34" = "produto"
    "D34" = "7"
    "E34" = "1"
    "F34" = "-6"
    "G34" = "puta"
    "H34" = "quatro"
    "I34" = "naotem"
    "J34" = "naotem"
}
foreach ($addr in $row34.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $row34[$addr]
    $cell.ClearFormats()
}
